# Add a "Type" column (column E) to the TEST_SET sheet, and fill in the
# one "Not Missing" value that row 3 already carries in columns A and D.
# This mirrors the commit: "added type column; made 4 owned col not required"
# (the "4 Owned" / required-column semantics live in the consuming Python
# code, not in the worksheet itself - the only workbook-visible change is
# the new Type column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1 = "Type" (becomes a new shared string).
$ws.Cells.Item(1, 5).Value = "Type"

# Row 3 already marks "Not Missing" in A3/D3 - do the same for the new
# Type column so the row stays consistent.
$ws.Cells.Item(3, 5).Value = "Not Missing"

# Give column E (and D, which picks up the same treatment in the diff)
# the same best-fit-ish width treatment column A already has.
$ws.Range("D1:E1").ColumnWidth = 9.67

# Selection follows the last-edited / newly added cell, as in the diff
# (activeCell="E4" sqref="E4").
$ws.Range("E4").Select() | Out-Null
